# Updates cryptos list data (price + 1h volume change columns, and a row swap
# of Monero/Fetch.AI) to match latest scrape, per commit message:
# "Updated cryptos list on Tue Mar 19 07:28:10 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.355.11'
$ws.Range("E2").Value = '  -5.72%  '
# Row 3
$ws.Range("D3").Value = '3.341.70'
$ws.Range("E3").Value = '  -7.64%  '
# Row 4
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.35%  '
# Row 5
$ws.Range("D5").Value = "'184.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -8.31%  '
# Row 6
$ws.Range("D6").Value = "'523.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -9.40%  '
# Row 7
$ws.Range("D7").Value = "'0.591"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.82%  '
# Row 8
$ws.Range("D8").Value = '3.334.23'
$ws.Range("E8").Value = '  -7.71%  '
# Row 9
$ws.Range("E9").Value = '  -0.01%  '
# Row 10
$ws.Range("D10").Value = "'0.614"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -10.17%  '
# Row 11
$ws.Range("D11").Value = "'56.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.72%  '
# Row 12
$ws.Range("D12").Value = "'0.130"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -12.93%  '
# Row 13
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -12.65%  '
# Row 14
$ws.Range("D14").Value = "'9.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -10.72%  '
# Row 15
$ws.Range("D15").Value = '3.874.36'
$ws.Range("E15").Value = '  -7.68%  '
# Row 16
$ws.Range("E16").Value = '  -4.54%  '
# Row 17
$ws.Range("D17").Value = '3.336.23'
$ws.Range("E17").Value = '  -8.00%  '
# Row 18
$ws.Range("D18").Value = '63.968.36'
$ws.Range("E18").Value = '  -6.09%  '
# Row 19
$ws.Range("D19").Value = "'17.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -11.12%  '
# Row 20
$ws.Range("D20").Value = "'10.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -12.67%  '
# Row 21
$ws.Range("D21").Value = "'0.951"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -11.49%  '
# Row 22
$ws.Range("D22").Value = "'371.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -8.58%  '
# Row 23
$ws.Range("D23").Value = "'79.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.81%  '
# Row 24
$ws.Range("D24").Value = "'3.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -13.77%  '
# Row 25
$ws.Range("D25").Value = "'10.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -17.93%  '
# Row 26
$ws.Range("D26").Value = "'5.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.46%  '
# Row 27
$ws.Range("D27").Value = "'3.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.33%  '
# Row 28
$ws.Range("D28").Value = "'2.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -11.07%  '
# Row 29
$ws.Range("D29").Value = "'11.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.49%  '
# Row 30
$ws.Range("D30").Value = "'8.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -11.71%  '
# Row 31
$ws.Range("D31").Value = "'653.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.98%  '
# Row 32
$ws.Range("D32").Value = "'28.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -10.08%  '
# Row 33
$ws.Range("D33").Value = "'6.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -14.81%  '
# Row 34
$ws.Range("D34").Value = "'10.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -10.22%  '
# Row 35
$ws.Range("D35").Value = "'58.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.01%  '
# Row 36
$ws.Range("D36").Value = "'0.102"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -10.63%  '
# Row 37
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.06%  '
# Row 38
$ws.Range("D38").Value = "'35.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -14.53%  '
# Row 39
$ws.Range("D39").Value = "'0.370"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -10.74%  '
# Row 40
$ws.Range("D40").Value = "'0.995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.26%  '
# Row 41
$ws.Range("D41").Value = "'0.123"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.11%  '
# Row 42
$ws.Range("D42").Value = "'2.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -15.31%  '
# Row 43
$ws.Range("D43").Value = '2.747.34'
$ws.Range("E43").Value = '  -14.10%  '
# Row 44
$ws.Range("E44").Value = '  -8.66%  '
# Row 45
$ws.Range("D45").Value = '0.0₃0608'
$ws.Range("E45").Value = '  -20.76%  '
# Row 46
$ws.Range("D46").Value = "'0.0379"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.08%  '
# Row 47
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = "'2.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -15.74%  '
# Row 48
$ws.Range("D48").Value = "'0.123"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.32%  '
# Row 49
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = "'134.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.32%  '
# Row 50
$ws.Range("D50").Value = "'2.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -10.89%  '
# Row 51
$ws.Range("D51").Value = "'2.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.74%  '
